# Update countries & provincias Spain
# - refresh the "datos actualizados" timestamp
# - update case counts for several countries
# - two pairs of countries swapped rank (their rows now carry the other's
#   former numbers alongside the new, updated figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 24 de Junio de 2020 a las 11:14'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2424493
$ws.Range("C4").Value = 325
$ws.Range("E4").Value = 1280605

# India (row 7)
$ws.Range("B7").Value = 456926
$ws.Range("C7").Value = 811
$ws.Range("D7").Value = 258827
$ws.Range("E7").Value = 183609
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 14490

# Banglades (row 20)
$ws.Range("B20").Value = 122660
$ws.Range("C20").Value = 3462
$ws.Range("D20").Value = 49666
$ws.Range("E20").Value = 71412
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 1582

# Polonia (row 40)
$ws.Range("B40").Value = 32821
$ws.Range("C40").Value = 294
$ws.Range("D40").Value = 18134
$ws.Range("E40").Value = 13291
$ws.Range("G40").Value = 21
$ws.Range("H40").Value = 1396

# Israel (row 52)
$ws.Range("B52").Value = 21666
$ws.Range("C52").Value = 154
$ws.Range("D52").Value = 15898
$ws.Range("E52").Value = 5460

# Austria (row 56)
$ws.Range("B56").Value = 17449
$ws.Range("C56").Value = 41
$ws.Range("D56").Value = 16282
$ws.Range("E56").Value = 474

# Sri Lanka (row 107)
$ws.Range("D107").Value = 1562
$ws.Range("E107").Value = 418

# Lituania (row 112)
$ws.Range("B112").Value = 1804
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 1484
$ws.Range("E112").Value = 242
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 78

# Eslovenia (row 118)
$ws.Range("B118").Value = 1541
$ws.Range("C118").Value = 7
$ws.Range("E118").Value = 56

# Hong Kong / Estado de Palestina swap rank (rows 123-124)
$ws.Range("A123").Value = 'Estado de Palestina'
$ws.Range("B123").Value = 1311
$ws.Range("C123").Value = 142
$ws.Range("D123").Value = 442
$ws.Range("E123").Value = 866
$ws.Range("H123").Value = 3

$ws.Range("A124").Value = 'Hong Kong'
$ws.Range("B124").Value = 1178
$ws.Range("D124").Value = 1083
$ws.Range("E124").Value = 89
$ws.Range("H124").Value = 6

# Cabo Verde (row 132)
$ws.Range("D132").Value = 479
$ws.Range("E132").Value = 495

# Georgia / Burkina Faso swap rank (rows 133-134)
$ws.Range("A133").Value = 'Burkina Faso'
$ws.Range("B133").Value = 919
$ws.Range("C133").Value = 12
$ws.Range("D133").Value = 825
$ws.Range("E133").Value = 41
$ws.Range("H133").Value = 53

$ws.Range("A134").Value = 'Georgia'
$ws.Range("B134").Value = 914
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 771
$ws.Range("E134").Value = 129
$ws.Range("H134").Value = 14

# Namibia (row 185)
$ws.Range("B185").Value = 76
$ws.Range("C185").Value = 4
$ws.Range("E185").Value = 55

# Fiyi / Dominica swap rank (rows 202-203), stats identical so only labels move
$ws.Range("A202").Value = 'Dominica'
$ws.Range("A203").Value = 'Fiyi'

# Seychelles / Montserrat swap rank (rows 211-212)
$ws.Range("A211").Value = 'Montserrat'
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("A212").Value = 'Seychelles'
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
